$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3438530600655295
$ws.Range("D2").Value = 50.9346037455865
$ws.Range("C3").Value = 0.4825802292055428
$ws.Range("D3").Value = 48.95278704358631
$ws.Range("E3").Value = 0.330963389234032
$ws.Range("C4").Value = 0.6697017583665801
$ws.Range("D4").Value = 46.27962234128579
$ws.Range("E4").Value = 0.7773818945182209
$ws.Range("C5").Value = 0.9151453693242683
$ws.Range("D5").Value = 42.77328504189024
$ws.Range("E5").Value = 1.362940223517277
$ws.Range("C6").Value = 1.224914635269494
$ws.Range("D6").Value = 38.34800981410131
$ws.Range("E6").Value = 2.101961186558029
$ws.Range("C7").Value = 1.596051013944228
$ws.Range("D7").Value = 33.0460615473194
$ws.Range("E7").Value = 2.987386547110608
$ws.Range("C8").Value = 2.011526501780981
$ws.Range("D8").Value = 27.11069743536578
$ws.Range("E8").Value = 3.978592353806863
$ws.Range("C9").Value = 2.439016549222489
$ws.Range("D9").Value = 21.00369675762995
$ws.Range("E9").Value = 4.998461466988747
$ws.Range("C10").Value = 2.837881536945993
$ws.Range("D10").Value = 15.30562550443704
$ws.Range("E10").Value = 5.950039366271963
$ws.Range("C11").Value = 3.173379749207618
$ws.Range("D11").Value = 10.51279390069954
$ws.Range("E11").Value = 6.750442244096125
$ws.Range("C12").Value = 3.429012208198724
$ws.Range("D12").Value = 6.860901629398032
$ws.Range("E12").Value = 7.360308253403478
$ws.Range("C13").Value = 3.608050188682537
$ws.Range("D13").Value = 4.303216193914981
$ws.Range("E13").Value = 7.787441721129148
$ws.Range("C14").Value = 3.725594175089091
$ws.Range("D14").Value = 2.624016388107058
$ws.Range("E14").Value = 8.067868088699072
$ws.Range("C15").Value = 3.799339808000172
$ws.Range("D15").Value = 1.570507346520189
$ws.Range("E15").Value = 8.243804098644079
$ws.Range("C16").Value = 3.844246665071113
$ws.Range("D16").Value = 0.9289808169353263
$ws.Range("E16").Value = 8.350939029084749
$ws.Range("C17").Value = 3.87108503956683
$ws.Range("D17").Value = 0.5455754669965129
$ws.Range("E17").Value = 8.414967722524532
$ws.Range("C18").Value = 3.886942998941307
$ws.Range("D18").Value = 0.3190331902182707
$ws.Range("E18").Value = 8.452800282746498
$ws.Range("C19").Value = 3.896249335110733
$ws.Range("D19").Value = 0.1860855306550447
$ws.Range("E19").Value = 8.475002541893558
$ws.Range("C20").Value = 3.901688867334186
$ws.Range("D20").Value = 0.1083779274628423
$ws.Range("E20").Value = 8.487979711626657
$ws.Range("C21").Value = 3.904860756908391
$ws.Range("D21").Value = 0.06306521925991308
$ws.Range("E21").Value = 8.495546933896545
